# Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite
#
# The "CasesTab" Cypher query (cell B2 on the "startup" sheet) had its last
# RETURN line - the `co.cohort_description` / `Cohort` column - removed, and
# the trailing comma on the preceding line ("Response to Treatment") dropped
# so the query still parses correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nMATCH (samp:sample)-->(c) `n  MATCH (f:file)-[*]->(c)`n    WHERE f.file_format IN [`"bam`"]  `nOPTIONAL MATCH (co:cohort)<-[*]-(c)`n  WITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

# The row is now one line shorter, so the workbook was re-saved with the
# view zoomed in further and the selection moved onto the edited cell, with
# the (now shorter) wrapped-text row heights re-fit to the new content.
$ws.Select()
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 85

$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 259.2
